# Auto-generated edit script: apply the "Updated cryptos list" diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.937.91"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.226.71"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.90"
$ws.Range("E5").Value = "  +3.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.75"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  +9.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.11"
$ws.Range("E10").Value = "  +11.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0972"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.32"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.29"
$ws.Range("E13").Value = "  +8.94%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "2.558.34"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.07"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "2.208.05"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").Value = "41.927.25"
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.99"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.62"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +13.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("E26").Value = "  +21.48%  "
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  -1.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.15"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.89"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.59"
$ws.Range("E34").Value = "  +6.29%  "
$ws.Range("E35").Value = "  +4.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.70"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.24"
$ws.Range("E37").Value = "  +15.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.13"
$ws.Range("E38").Value = "  +9.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0313"
$ws.Range("E39").Value = "  +10.18%  "
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.46"
$ws.Range("E42").Value = "  +24.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.89"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.205"
$ws.Range("E44").Value = "  +8.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.89"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.90"
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.64"
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +6.44%  "
$ws.Range("E51").Value = "  +5.31%  "
